# Excel COM-interop script implementing the AssetList.xlsx edit:
#  1. Mark the "LowHealth" row's Status as "Completed" (was "In progress").
#  2. Add two new rows documenting the UI click/hover sounds:
#       UIClick / Sound of clicking on a button / Interface / Click SFX (x2) / In progress
#       UIHover / Sound of hovering over a button / Interface / Hover SFX (x2) / In progress
#  3. Leave the selection on F27 (the last cell touched), matching the
#     author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. LowHealth status -> Completed -------------------------------------
$ws.Range("E24").Value = "Completed"

# --- 2. Append the two new asset rows --------------------------------------
# Values are written column-by-column (not row-by-row) so new shared-string
# entries land in the same order the source workbook used.
$ws.Range("A26").Value = "UIClick"
$ws.Range("A27").Value = "UIHover"

$ws.Range("B26").Value = "Sound of clicking on a button"
$ws.Range("B27").Value = "Sound of hovering over a button"

$ws.Range("D26").Value = "Click SFX (x2)"
$ws.Range("D27").Value = "Hover SFX (x2)"

$ws.Range("C26").Value = "Interface"
$ws.Range("C27").Value = "Interface"

$ws.Range("E26").Value = "In progress"
$ws.Range("E27").Value = "In progress"

# --- 3. Update the visible selection ---------------------------------------
$ws.Range("F27").Select() | Out-Null
